$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 (Objetivos:) -- B/C now hold the responsible professor entry
$ws.Range("B10").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C10").Value = "519033 - Carlos Yujiro Shigue"

# Row 13 (Programa resumido:) -- B/C now hold the activation date, copied
# from B8 (which already stores "01/01/2023" as literal text) so it is not
# re-interpreted as a date value/format by Excel's text-to-date heuristics.
$ws.Range("B8").Copy()
$ws.Range("B13").PasteSpecial(-4163)
$ws.Range("B8").Copy()
$ws.Range("C13").PasteSpecial(-4163)

# Row 15 (Programa:) -- B/C now hold the same professor entry as row 10
$ws.Range("B15").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C15").Value = "519033 - Carlos Yujiro Shigue"

# Row 18 (Método:) -- B/C now hold the second professor entry
$ws.Range("B18").Value = "7290967 - Emerson Gonçalves de Melo"
$ws.Range("C18").Value = "7290967 - Emerson Gonçalves de Melo"

$excel.CutCopyMode = $false
